# Daily attendance processing - 2025-10-17 02:48:58
# Normalize the "Recorded By" (column G) entries so that "System" is
# listed first among the comma-separated recorder names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2   = "System, backup@backdoor.com, system"
    4   = "System, backup@backdoor.com"
    5   = "System, backup@backdoor.com"
    10  = "System, dnasr281@gmail.com"
    18  = "System, dnasr281@gmail.com"
    19  = "System, dnasr281@gmail.com"
    29  = "System, backup@backdoor.com, system"
    31  = "System, backup@backdoor.com"
    32  = "System, backup@backdoor.com"
    37  = "System, dnasr281@gmail.com"
    45  = "System, dnasr281@gmail.com"
    46  = "System, dnasr281@gmail.com"
    56  = "System, backup@backdoor.com, system"
    58  = "System, backup@backdoor.com"
    59  = "System, backup@backdoor.com"
    64  = "System, dnasr281@gmail.com"
    72  = "System, dnasr281@gmail.com"
    73  = "System, dnasr281@gmail.com"
    83  = "System, backup@backdoor.com"
    84  = "System, backup@backdoor.com"
    85  = "System, backup@backdoor.com"
    86  = "System, dnasr281@gmail.com"
    97  = "System, dnasr281@gmail.com"
    109 = "System, backup@backdoor.com"
    110 = "System, backup@backdoor.com"
    111 = "System, backup@backdoor.com"
    112 = "System, dnasr281@gmail.com"
    123 = "System, dnasr281@gmail.com"
    135 = "System, backup@backdoor.com"
    136 = "System, backup@backdoor.com"
    137 = "System, backup@backdoor.com"
    138 = "System, dnasr281@gmail.com"
    149 = "System, dnasr281@gmail.com"
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
